$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.204.35"
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("D3").Value = "2.385.45"
$ws.Range("E3").Value = "  -2.57%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.12%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  +0.60%  "

$ws.Range("D9").Value = "2.386.91"

$ws.Range("E10").Value = "  -4.30%  "

$ws.Range("E11").Value = "  -1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("E13").Value = "  -1.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "

$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("E16").Value = "  -3.57%  "

$ws.Range("D17").Value = "60.231.68"
$ws.Range("E17").Value = "  -2.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +15.63%  "

$ws.Range("D19").Value = "2.405.34"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "552.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.97%  "

$ws.Range("D29").Value = "2.508.81"
$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("D30").Value = "0.0₃0900"
$ws.Range("E30").Value = "  -2.69%  "

$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  -5.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.01%  "

$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.85%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("E38").Value = "  -1.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.75%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("E44").Value = "  -2.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.24%  "

$ws.Range("E46").Value = "  -6.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("E48").Value = "  -2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.587"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.38%  "
